$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the "Manufacturer Part Number" (column G) values for the rows that
# correspond to designators J2, J3, R19, R20, SW1, U2, U3 (rows 8, 9, 13, 14, 15, 17, 18).
# Clearing contents also resets those cells back to the same style used by the
# rest of column G (style used by e.g. G10-G12).
$rows = @(8, 9, 13, 14, 15, 17, 18)
foreach ($r in $rows) {
    $cell = $ws.Cells.Item($r, 7)
    $cell.Value = ""
}

# A handful of these cells (G8, G9, G15, G17) previously used a different,
# one-off font/format (no border, 9pt Segoe UI). Re-apply the standard
# formatting used by the rest of column G (8pt Segoe UI, black, thin border)
# by copying the format from a neighboring cell that already uses it.
$fmtRows = @(8, 9, 15, 17)
$source = $ws.Cells.Item(10, 7)
$source.Copy() | Out-Null
foreach ($r in $fmtRows) {
    $target = $ws.Cells.Item($r, 7)
    $target.PasteSpecial(-4122)  # xlPasteFormats
}
$excel.CutCopyMode = 0

